# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps for the first (row 2) entry on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-12 22:57:11"
$wsZhCn.Range("H2").Value = "2016-03-12 22:57:27"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-12 22:57:14"
$wsDeDe.Range("H2").Value = "2016-03-12 22:57:33"
